$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.601.14"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.142.16"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.99"
$ws.Range("E5").Value = "  +5.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5265"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4578"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.55"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09232"
$ws.Range("E10").Value = "  +3.78%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.57"
$ws.Range("E12").Value = "  +6.32%  "
$ws.Range("D13").Value = "2.133.41"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.886"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.176"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.07"
$ws.Range("E16").Value = "  +5.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001175"
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06717"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.94"
$ws.Range("E20").Value = "  +4.06%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.365"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").Value = "30.700.72"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.94"
$ws.Range("E24").Value = "  +4.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.386"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").Value = "2.371.51"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.60"
$ws.Range("E27").Value = "  +2.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.663"
$ws.Range("E28").Value = "  +5.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.74"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "137.44"
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1084"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.672"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.379"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.026"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("E36").Value = "  +6.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.48"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02654"
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07018"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2346"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.70"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7018"
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.81"
$ws.Range("E44").Value = "  +6.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.368"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6506"
$ws.Range("E46").Value = "  +2.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.753"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000369"
$ws.Range("E48").Value = "  +5.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.256"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.32"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07308"
$ws.Range("E51").Value = "  +2.58%  "
